$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by 10 days
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 10
}

# Updated Actual Production (MW) readings for the new date range (rows 2-39)
$newValues = @(876, 834, 928, 1037, 1166, 1181, 1128, 1071, 953, 853, 786, 803, 832, 868, 926, 938, 1001, 1064, 1104, 1172, 1279, 1349, 1379, 1419, 1494, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

$r = 2
foreach ($val in $newValues) {
    $ws.Cells.Item($r, 2).Value2 = $val
    $r++
}
